# Add a "Save" column (H) to the s_vals sheet, mirroring the existing
# header style used by the other column headers (B1:G1) and filling in
# the per-row save indicator values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - same bold/bordered/centered style as the other headers.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Per-row "Save" values for H2:H16.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}
